# Update the Angpt1-Tie1 LR-pairs sheet with refreshed TPM-derived NATMI output.
# All data rows (2-16) are replaced: sending/target cluster labels were
# reshuffled (ECs now appears as a sending cluster) and every numeric metric
# column (E:T) was recomputed from the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Angpt1"
$ws.Range("C2").Value = "Tie1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = [double]"0.6666666666666666"
$ws.Range("G2").Value = [double]"0.1497083333333333"
$ws.Range("H2").Value = [double]"0.449125"
$ws.Range("I2").Value = [double]"0.006513369349540601"
$ws.Range("J2").Value = [double]"0.006769619242096868"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = [double]"56.89751433333333"
$ws.Range("N2").Value = [double]"170.692543"
$ws.Range("O2").Value = [double]"0.9873538794860484"
$ws.Range("P2").Value = [double]"0.9895831360385335"
$ws.Range("Q2").Value = [double]"8.518032041652777"
$ws.Range("R2").Value = [double]"76.662288374875"
$ws.Range("S2").Value = [double]"0.006431000495794432"
$ws.Range("T2").Value = [double]"0.006699101039381019"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Angpt1"
$ws.Range("C3").Value = "Tie1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = [double]"0.6666666666666666"
$ws.Range("G3").Value = [double]"0.1497083333333333"
$ws.Range("H3").Value = [double]"0.449125"
$ws.Range("I3").Value = [double]"0.006513369349540601"
$ws.Range("J3").Value = [double]"0.006769619242096868"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = [double]"0.3333333333333333"
$ws.Range("M3").Value = [double]"0.1599326666666667"
$ws.Range("N3").Value = [double]"0.479798"
$ws.Range("O3").Value = [double]"0.002775343364997773"
$ws.Range("P3").Value = [double]"0.002781609560442347"
$ws.Range("Q3").Value = [double]"0.02394325297222222"
$ws.Range("R3").Value = [double]"0.21548927675"
$ws.Range("S3").Value = [double]"1.807683640802737e-05"
$ws.Range("T3").Value = [double]"1.883043760437112e-05"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Angpt1"
$ws.Range("C4").Value = "Tie1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = [double]"0.6666666666666666"
$ws.Range("G4").Value = [double]"0.1497083333333333"
$ws.Range("H4").Value = [double]"0.449125"
$ws.Range("I4").Value = [double]"0.006513369349540601"
$ws.Range("J4").Value = [double]"0.006769619242096868"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = [double]"0.3333333333333333"
$ws.Range("M4").Value = [double]"0.1362956666666667"
$ws.Range("N4").Value = [double]"0.408887"
$ws.Range("O4").Value = [double]"0.002365165804117242"
$ws.Range("P4").Value = [double]"0.002370505896941191"
$ws.Range("Q4").Value = [double]"0.02040459709722222"
$ws.Range("R4").Value = [double]"0.183641373875"
$ws.Range("S4").Value = [double]"1.540519845511879e-05"
$ws.Range("T4").Value = [double]"1.604742233343718e-05"

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Angpt1"
$ws.Range("C5").Value = "Tie1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = [double]"0.6666666666666666"
$ws.Range("G5").Value = [double]"0.1497083333333333"
$ws.Range("H5").Value = [double]"0.449125"
$ws.Range("I5").Value = [double]"0.006513369349540601"
$ws.Range("J5").Value = [double]"0.006769619242096868"
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = [double]"0.389448"
$ws.Range("N5").Value = [double]"0.778896"
$ws.Range("O5").Value = [double]"0.006758168580183656"
$ws.Range("P5").Value = [double]"0.004515618156370602"
$ws.Range("Q5").Value = [double]"0.05830361100000001"
$ws.Range("R5").Value = [double]"0.349821666"
$ws.Range("S5").Value = [double]"4.401844808919655e-05"
$ws.Range("T5").Value = [double]"3.056901556132841e-05"

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Angpt1"
$ws.Range("C6").Value = "Tie1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = [double]"0.6666666666666666"
$ws.Range("G6").Value = [double]"0.1497083333333333"
$ws.Range("H6").Value = [double]"0.449125"
$ws.Range("I6").Value = [double]"0.006513369349540601"
$ws.Range("J6").Value = [double]"0.006769619242096868"
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.04307233333333333"
$ws.Range("N6").Value = [double]"0.129217"
$ws.Range("O6").Value = [double]"0.0007474427646528689"
$ws.Range("P6").Value = [double]"0.0007491303477123261"
$ws.Range("Q6").Value = [double]"0.00644828723611111"
$ws.Range("R6").Value = [double]"0.058034585125"
$ws.Range("S6").Value = [double]"4.868370793825885e-06"
$ws.Range("T6").Value = [double]"5.071327216712081e-06"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Angpt1"
$ws.Range("C7").Value = "Tie1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = [double]"20.22494433333334"
$ws.Range("H7").Value = [double]"60.67483300000001"
$ws.Range("I7").Value = [double]"0.8799278542737426"
$ws.Range("J7").Value = [double]"0.9145460996110527"
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = [double]"56.89751433333333"
$ws.Range("N7").Value = [double]"170.692543"
$ws.Range("O7").Value = [double]"0.9873538794860484"
$ws.Range("P7").Value = [double]"0.9895831360385335"
$ws.Range("Q7").Value = [double]"1150.749060096702"
$ws.Range("R7").Value = [double]"10356.74154087032"
$ws.Range("S7").Value = [double]"0.8688001805850141"
$ws.Range("T7").Value = [double]"0.9050193973049147"

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Angpt1"
$ws.Range("C8").Value = "Tie1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = [double]"20.22494433333334"
$ws.Range("H8").Value = [double]"60.67483300000001"
$ws.Range("I8").Value = [double]"0.8799278542737426"
$ws.Range("J8").Value = [double]"0.9145460996110527"
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = [double]"0.3333333333333333"
$ws.Range("M8").Value = [double]"0.1599326666666667"
$ws.Range("N8").Value = [double]"0.479798"
$ws.Range("O8").Value = [double]"0.002775343364997773"
$ws.Range("P8").Value = [double]"0.002781609560442347"
$ws.Range("Q8").Value = [double]"3.23462928041489"
$ws.Range("R8").Value = [double]"29.111663523734"
$ws.Range("S8").Value = [double]"0.002442101932035359"
$ws.Range("T8").Value = [double]"0.002543910174143363"

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Angpt1"
$ws.Range("C9").Value = "Tie1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = [double]"20.22494433333334"
$ws.Range("H9").Value = [double]"60.67483300000001"
$ws.Range("I9").Value = [double]"0.8799278542737426"
$ws.Range("J9").Value = [double]"0.9145460996110527"
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.1362956666666667"
$ws.Range("N9").Value = [double]"0.408887"
$ws.Range("O9").Value = [double]"0.002365165804117242"
$ws.Range("P9").Value = [double]"0.002370505896941191"
$ws.Range("Q9").Value = [double]"2.75657227120789"
$ws.Range("R9").Value = [double]"24.809150440871"
$ws.Range("S9").Value = [double]"0.002081175271018516"
$ws.Range("T9").Value = [double]"0.002167936922152567"

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Angpt1"
$ws.Range("C10").Value = "Tie1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = [double]"20.22494433333334"
$ws.Range("H10").Value = [double]"60.67483300000001"
$ws.Range("I10").Value = [double]"0.8799278542737426"
$ws.Range("J10").Value = [double]"0.9145460996110527"
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = [double]"0.389448"
$ws.Range("N10").Value = [double]"0.778896"
$ws.Range("O10").Value = [double]"0.006758168580183656"
$ws.Range("P10").Value = [double]"0.004515618156370602"
$ws.Range("Q10").Value = [double]"7.876564120728002"
$ws.Range("R10").Value = [double]"47.259384724368"
$ws.Range("S10").Value = [double]"0.00594670077758123"
$ws.Range("T10").Value = [double]"0.004129740972241587"

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Angpt1"
$ws.Range("C11").Value = "Tie1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = [double]"20.22494433333334"
$ws.Range("H11").Value = [double]"60.67483300000001"
$ws.Range("I11").Value = [double]"0.8799278542737426"
$ws.Range("J11").Value = [double]"0.9145460996110527"
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = [double]"0.6666666666666666"
$ws.Range("M11").Value = [double]"0.04307233333333333"
$ws.Range("N11").Value = [double]"0.129217"
$ws.Range("O11").Value = [double]"0.0007474427646528689"
$ws.Range("P11").Value = [double]"0.0007491303477123261"
$ws.Range("Q11").Value = [double]"0.8711355439734445"
$ws.Range("R11").Value = [double]"7.840219895761001"
$ws.Range("S11").Value = [double]"0.0006576957080934329"
$ws.Range("T11").Value = [double]"0.0006851142376005795"

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Angpt1"
$ws.Range("C12").Value = "Tie1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = [double]"2.610123"
$ws.Range("H12").Value = [double]"5.220245999999999"
$ws.Range("I12").Value = [double]"0.1135587763767167"
$ws.Range("J12").Value = [double]"0.07868428114685043"
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = [double]"56.89751433333333"
$ws.Range("N12").Value = [double]"170.692543"
$ws.Range("O12").Value = [double]"0.9873538794860484"
$ws.Range("P12").Value = [double]"0.9895831360385335"
$ws.Range("Q12").Value = [double]"148.509510804263"
$ws.Range("R12").Value = [double]"891.0570648255779"
$ws.Range("S12").Value = [double]"0.1121226984052398"
$ws.Range("T12").Value = [double]"0.0778646376942379"

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Angpt1"
$ws.Range("C13").Value = "Tie1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = [double]"2.610123"
$ws.Range("H13").Value = [double]"5.220245999999999"
$ws.Range("I13").Value = [double]"0.1135587763767167"
$ws.Range("J13").Value = [double]"0.07868428114685043"
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.1599326666666667"
$ws.Range("N13").Value = [double]"0.479798"
$ws.Range("O13").Value = [double]"0.002775343364997773"
$ws.Range("P13").Value = [double]"0.002781609560442347"
$ws.Range("Q13").Value = [double]"0.417443931718"
$ws.Range("R13").Value = [double]"2.504663590308"
$ws.Range("S13").Value = [double]"0.0003151645965543865"
$ws.Range("T13").Value = [double]"0.0002188689486946126"

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Angpt1"
$ws.Range("C14").Value = "Tie1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = [double]"2.610123"
$ws.Range("H14").Value = [double]"5.220245999999999"
$ws.Range("I14").Value = [double]"0.1135587763767167"
$ws.Range("J14").Value = [double]"0.07868428114685043"
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = [double]"0.3333333333333333"
$ws.Range("M14").Value = [double]"0.1362956666666667"
$ws.Range("N14").Value = [double]"0.408887"
$ws.Range("O14").Value = [double]"0.002365165804117242"
$ws.Range("P14").Value = [double]"0.002370505896941191"
$ws.Range("Q14").Value = [double]"0.355748454367"
$ws.Range("R14").Value = [double]"2.134490726202"
$ws.Range("S14").Value = [double]"0.0002685853346436072"
$ws.Range("T14").Value = [double]"0.0001865215524551875"

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Angpt1"
$ws.Range("C15").Value = "Tie1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = [double]"2.610123"
$ws.Range("H15").Value = [double]"5.220245999999999"
$ws.Range("I15").Value = [double]"0.1135587763767167"
$ws.Range("J15").Value = [double]"0.07868428114685043"
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = [double]"0.389448"
$ws.Range("N15").Value = [double]"0.778896"
$ws.Range("O15").Value = [double]"0.006758168580183656"
$ws.Range("P15").Value = [double]"0.004515618156370602"
$ws.Range("Q15").Value = [double]"1.016507182104"
$ws.Range("R15").Value = [double]"4.066028728416"
$ws.Range("S15").Value = [double]"0.0007674493545132288"
$ws.Range("T15").Value = [double]"0.0003553081685676869"

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Angpt1"
$ws.Range("C16").Value = "Tie1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = [double]"2.610123"
$ws.Range("H16").Value = [double]"5.220245999999999"
$ws.Range("I16").Value = [double]"0.1135587763767167"
$ws.Range("J16").Value = [double]"0.07868428114685043"
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = [double]"0.6666666666666666"
$ws.Range("M16").Value = [double]"0.04307233333333333"
$ws.Range("N16").Value = [double]"0.129217"
$ws.Range("O16").Value = [double]"0.0007474427646528689"
$ws.Range("P16").Value = [double]"0.0007491303477123261"
$ws.Range("Q16").Value = [double]"0.112424087897"
$ws.Range("R16").Value = [double]"0.6745445273819999"
$ws.Range("S16").Value = [double]"8.487868576561001e-05"
$ws.Range("T16").Value = [double]"5.894478289503449e-05"
